# Replace the three-digit division expressions in the document's table
# according to the commit diff. Each old expression is unique in the
# document, so a straightforward literal Find & Replace (ReplaceAll) is
# safe for each pair.

$d = $word.ActiveDocument

$replacements = @(
    @("282÷8=", "836÷6="),
    @("543÷7=", "331÷8="),
    @("548÷2=", "453÷4="),
    @("806÷2=", "619÷6="),
    @("431÷2=", "340÷5="),
    @("216÷6=", "685÷9="),
    @("161÷4=", "739÷5="),
    @("651÷9=", "375÷3="),
    @("908÷8=", "110÷5="),
    @("262÷4=", "761÷4="),
    @("253÷8=", "692÷7="),
    @("309÷7=", "829÷2="),
    @("649÷2=", "178÷7="),
    @("848÷2=", "981÷5="),
    @("279÷9=", "734÷5="),
    @("253÷7=", "682÷5="),
    @("936÷6=", "216÷9="),
    @("508÷7=", "247÷8="),
    @("356÷7=", "509÷4="),
    @("342÷8=", "593÷7="),
    @("471÷5=", "461÷6="),
    @("257÷8=", "850÷3="),
    @("639÷4=", "130÷4="),
    @("567÷5=", "588÷9="),
    @("382÷8=", "273÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
